$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fold the old A2/A3/A4 rows into a single A2 value, then drop the now-empty rows 3-4
$ws.Range("A2").Value = "('Elemental Shaman', ['Token Creature — Elemental Shaman', '3/1'])"
$ws.Rows("3:4").Delete()
